# Update weekly price data for "Hortaliza, Terminal Hortofrutícola Agro Chillán - Alcachofa"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value2 = 44432
$ws.Range("O2").Value2 = 'Provincia del Elquí'

# Row 3
$ws.Range("D3").Value2 = 44454
$ws.Range("J3").Value2 = 120
$ws.Range("K3").Value2 = 13000
$ws.Range("L3").Value2 = 14000
$ws.Range("M3").Value2 = 13500
$ws.Range("P3").Value2 = 338

# Row 4
$ws.Range("D4").Value2 = 44435
$ws.Range("K4").Value2 = 14000
$ws.Range("L4").Value2 = 15000
$ws.Range("M4").Value2 = 14500
$ws.Range("O4").Value2 = 'Provincia del Elquí'
$ws.Range("P4").Value2 = 362

# Row 5
$ws.Range("D5").Value2 = 44446
$ws.Range("J5").Value2 = 160
$ws.Range("K5").Value2 = 12500
$ws.Range("L5").Value2 = 13000
$ws.Range("M5").Value2 = 12750
$ws.Range("P5").Value2 = 319

# Row 6
$ws.Range("D6").Value2 = 44417
$ws.Range("K6").Value2 = 15000
$ws.Range("L6").Value2 = 16000
$ws.Range("M6").Value2 = 15500
$ws.Range("P6").Value2 = 388

# Row 7
$ws.Range("D7").Value2 = 44399
$ws.Range("H7").Value2 = 'Española'
$ws.Range("I7").Value2 = 'Segunda'
$ws.Range("K7").Value2 = 15500
$ws.Range("M7").Value2 = 15750
$ws.Range("P7").Value2 = 394

# Row 8
$ws.Range("D8").Value2 = 44426
$ws.Range("K8").Value2 = 13000
$ws.Range("L8").Value2 = 14000
$ws.Range("M8").Value2 = 13500
$ws.Range("O8").Value2 = 'Región del Maule'
$ws.Range("P8").Value2 = 338

# Row 9
$ws.Range("D9").Value2 = 44453

# Row 10
$ws.Range("D10").Value2 = 44420
$ws.Range("K10").Value2 = 13000
$ws.Range("L10").Value2 = 14000
$ws.Range("M10").Value2 = 13500
$ws.Range("P10").Value2 = 338

# Row 11
$ws.Range("D11").Value2 = 44427
$ws.Range("H11").Value2 = 'Madrigal'
$ws.Range("I11").Value2 = 'Primera'
$ws.Range("K11").Value2 = 13000
$ws.Range("L11").Value2 = 14000
$ws.Range("M11").Value2 = 13500
$ws.Range("P11").Value2 = 338

# New row 12 (duplicate of the original row 2 content, before today's edits)
$ws.Range("A12").Value2 = 7
$ws.Range("B12").Value2 = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C12").Value2 = 'Ñuble'
$ws.Range("D12").Value2 = 44425
$ws.Range("D12").NumberFormat = $ws.Range("D11").NumberFormat
$ws.Range("E12").Value2 = 16
$ws.Range("F12").Value2 = 100112013
$ws.Range("G12").Value2 = 'Alcachofa'
$ws.Range("H12").Value2 = 'Madrigal'
$ws.Range("I12").Value2 = 'Primera'
$ws.Range("J12").Value2 = 120
$ws.Range("K12").Value2 = 14000
$ws.Range("L12").Value2 = 15000
$ws.Range("M12").Value2 = 14500
$ws.Range("N12").Value2 = '$/caja 40 unidades'
$ws.Range("O12").Value2 = 'Región del Maule'
$ws.Range("P12").Value2 = 362
$ws.Range("Q12").Value2 = 40
$ws.Range("R12").Value2 = 'Hortaliza'
